# Auto-generated edit script applying the diff's cell-value changes
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4555
$ws.Range("I32").Value = 4555
$ws.Range("K32").Value = 4555
$ws.Range("M32").Value = -4229
$ws.Range("H116").Value = 9333.333000000001
$ws.Range("I116").Value = 9333.333000000001
$ws.Range("K116").Value = 9333.333000000001
$ws.Range("M116").Value = -5891.333000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 168.83333
$ws.Range("I5").Value = 177.94737
$ws.Range("J5").Value = 134.2
$ws.Range("K5").Value = 177.94737
$ws.Range("L5").Value = 134.2
$ws.Range("M5").Value = -65.94737000000001
$ws.Range("N5").Value = -358.2
$ws.Range("H50").Value = 29482.666
$ws.Range("I50").Value = 448
$ws.Range("K50").Value = 448
$ws.Range("M50").Value = 266
$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 2000
$ws.Range("K61").Value = 2000
$ws.Range("M61").Value = -1788
$ws.Range("H102").Value = 15001763
$ws.Range("I102").Value = 1001488.8
$ws.Range("K102").Value = 1001488.8
$ws.Range("M102").Value = -999866.8
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = ""
$ws.Range("H122").Value = 1374
$ws.Range("I122").Value = 1374
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4122
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1672
$ws.Range("N122").Value = ""
$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450
$ws.Range("I140").Value = 50429
$ws.Range("K140").Value = 50429
$ws.Range("M140").Value = -45249

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 168.83333
$ws.Range("I4").Value = 177.94737
$ws.Range("J4").Value = 134.2
$ws.Range("K4").Value = 177.94737
$ws.Range("L4").Value = 134.2
$ws.Range("M4").Value = -62.94737000000001
$ws.Range("N4").Value = -364.2
$ws.Range("H5").Value = 1716.2
$ws.Range("I5").Value = 300
$ws.Range("J5").Value = 2070.25
$ws.Range("K5").Value = 300
$ws.Range("L5").Value = 2070.25
$ws.Range("M5").Value = -187
$ws.Range("N5").Value = -2296.25
$ws.Range("H7").Value = 920
$ws.Range("I7").Value = 60
$ws.Range("K7").Value = 60
$ws.Range("M7").Value = 53
$ws.Range("H86").Value = 3614.2856
$ws.Range("I86").Value = 2900
$ws.Range("K86").Value = 2900
$ws.Range("M86").Value = -1777
$ws.Range("H89").Value = 3614.2856
$ws.Range("I89").Value = 2900
$ws.Range("K89").Value = 14500
$ws.Range("M89").Value = -8884

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 822.6
$ws.Range("J2").Value = 505
$ws.Range("L2").Value = 505
$ws.Range("N2").Value = -731
$ws.Range("H3").Value = 1499.75
$ws.Range("I3").Value = 1333
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 1333
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -1220
$ws.Range("N3").Value = -2226
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = ""
$ws.Range("H12").Value = 145.8
$ws.Range("I12").Value = 145.8
$ws.Range("K12").Value = 145.8
$ws.Range("M12").Value = 24.19999999999999
$ws.Range("H14").Value = 799
$ws.Range("J14").Value = 799
$ws.Range("L14").Value = 799
$ws.Range("N14").Value = -1139
$ws.Range("H15").Value = 7157.8
$ws.Range("I15").Value = 1790
$ws.Range("J15").Value = 8499.75
$ws.Range("K15").Value = 1790
$ws.Range("L15").Value = 8499.75
$ws.Range("M15").Value = -1620
$ws.Range("N15").Value = -8839.75
$ws.Range("H17").Value = 400
$ws.Range("I17").Value = 400
$ws.Range("K17").Value = 400
$ws.Range("M17").Value = -226
$ws.Range("H19").Value = 8571717
$ws.Range("I19").Value = 8571717
$ws.Range("K19").Value = 8571717
$ws.Range("M19").Value = -8571547
$ws.Range("H23").Value = 4000000
$ws.Range("I23").Value = 4000000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 4000000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -3999760
$ws.Range("N23").Value = ""
$ws.Range("H24").Value = 8571717
$ws.Range("I24").Value = 8571717
$ws.Range("K24").Value = 8571717
$ws.Range("M24").Value = -8571547
$ws.Range("H27").Value = 4000000
$ws.Range("I27").Value = 4000000
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 4000000
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -3999808
$ws.Range("N27").Value = ""
$ws.Range("H38").Value = 36000
$ws.Range("J38").Value = 37000
$ws.Range("L38").Value = 37000
$ws.Range("N38").Value = -37754
$ws.Range("H46").Value = 36000
$ws.Range("J46").Value = 37000
$ws.Range("L46").Value = 37000
$ws.Range("N46").Value = -37422
$ws.Range("H50").Value = 34500
$ws.Range("I50").Value = 20000
$ws.Range("J50").Value = 49000
$ws.Range("K50").Value = 20000
$ws.Range("L50").Value = 49000
$ws.Range("M50").Value = -19375
$ws.Range("N50").Value = -50250
$ws.Range("H51").Value = 45000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 45000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 45000
$ws.Range("M51").Value = ""
$ws.Range("N51").Value = -46472
$ws.Range("H61").Value = 45000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 45000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 45000
$ws.Range("M61").Value = ""
$ws.Range("N61").Value = -45696
$ws.Range("H122").Value = 3210.1538
$ws.Range("I122").Value = 914.7778
$ws.Range("K122").Value = 2744.3334
$ws.Range("M122").Value = -294.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 107
$ws.Range("I23").Value = 37.714287
$ws.Range("J23").Value = 151.09091
$ws.Range("K23").Value = 113.142861
$ws.Range("L23").Value = 453.27273
$ws.Range("M23").Value = 121.857139
$ws.Range("N23").Value = -923.27273
$ws.Range("H34").Value = 2433.3
$ws.Range("J34").Value = 2460.966
$ws.Range("L34").Value = 7382.897999999999
$ws.Range("N34").Value = -7550.897999999999
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = ""
$ws.Range("H131").Value = 2239.2856
$ws.Range("I131").Value = 1024
$ws.Range("J131").Value = 2914.4443
$ws.Range("K131").Value = 3072
$ws.Range("L131").Value = 8743.332900000001
$ws.Range("M131").Value = 1968
$ws.Range("N131").Value = -18823.3329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 32000
$ws.Range("J39").Value = 32000
$ws.Range("L39").Value = 32000
$ws.Range("N39").Value = -33064

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 33335930
$ws.Range("I93").Value = 55558292
$ws.Range("J93").Value = 2388.75
$ws.Range("K93").Value = 55558292
$ws.Range("L93").Value = 2388.75
$ws.Range("M93").Value = -55557044
$ws.Range("N93").Value = -4884.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492
